$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Fstl1", "Dip2a", "ECs", 3, 1, 21.30703333333333, 63.9211, 0.03229214412626456, 0.03229214412626456, 3, 1, 9.971009, 29.913027, 0.3700016197528959, 0.3700016197528959, 212.4526211299667, 1912.0735901697, 0.01194814563201185, 0.01194814563201185),
    @("ECs", "Fstl1", "Dip2a", "FAPs", 3, 1, 21.30703333333333, 63.9211, 0.03229214412626456, 0.03229214412626456, 3, 1, 10.04188433333333, 30.125653, 0.3726316432674529, 0.3726316432674529, 213.9627642198111, 1925.6648779783, 0.01203307473039939, 0.01203307473039939),
    @("ECs", "Fstl1", "Dip2a", "M2", 3, 1, 21.30703333333333, 63.9211, 0.03229214412626456, 0.03229214412626456, 3, 1, 2.166912666666667, 6.500738, 0.0804092340634467, 0.0804092340634467, 46.17048041908889, 415.5343237718, 0.002596586575459363, 0.002596586575459363),
    @("ECs", "Fstl1", "Dip2a", "sCs", 3, 1, 21.30703333333333, 63.9211, 0.03229214412626456, 0.03229214412626456, 3, 1, 4.768749000000001, 14.306247, 0.1769575029162046, 0.1769575029162046, 101.6078939013, 914.4710451117001, 0.005714337188393962, 0.005714337188393962),
    @("FAPs", "Fstl1", "Dip2a", "ECs", 3, 1, 614.8671876666667, 1844.601563, 0.9318697508010483, 0.9318697508010483, 3, 1, 9.971009, 29.913027, 0.3700016197528959, 0.3700016197528959, 6130.846262029023, 55177.61635826121, 0.3447933171951153, 0.3447933171951153),
    @("FAPs", "Fstl1", "Dip2a", "FAPs", 3, 1, 614.8671876666667, 1844.601563, 0.9318697508010483, 0.9318697508010483, 3, 1, 10.04188433333333, 30.125653, 0.3726316432674529, 0.3726316432674529, 6174.425178910627, 55569.82661019565, 0.3472441565522265, 0.3472441565522265),
    @("FAPs", "Fstl1", "Dip2a", "M2", 3, 1, 614.8671876666667, 1844.601563, 0.9318697508010483, 0.9318697508010483, 3, 1, 2.166912666666667, 6.500738, 0.0804092340634467, 0.0804092340634467, 1332.363497272611, 11991.2714754535, 0.07493093290880724, 0.07493093290880724),
    @("FAPs", "Fstl1", "Dip2a", "sCs", 3, 1, 614.8671876666667, 1844.601563, 0.9318697508010483, 0.9318697508010483, 3, 1, 4.768749000000001, 14.306247, 0.1769575029162046, 0.1769575029162046, 2932.14728631823, 26389.32557686407, 0.1649013441448994, 0.1649013441448994),
    @("M2", "Fstl1", "Dip2a", "ECs", 2, 0.6666666666666666, 0.327557, 0.982671, 0.0004964331583890222, 0.0004964331583890221, 3, 1, 9.971009, 29.913027, 0.3700016197528959, 0.3700016197528959, 3.266073795013, 29.394664155117, 0.0001836810727029841, 0.0001836810727029841),
    @("M2", "Fstl1", "Dip2a", "FAPs", 2, 0.6666666666666666, 0.327557, 0.982671, 0.0004964331583890222, 0.0004964331583890221, 3, 1, 10.04188433333333, 30.125653, 0.3726316432674529, 0.3726316432674529, 3.289289506573667, 29.603605559163, 0.0001849867035829531, 0.000184986703582953),
    @("M2", "Fstl1", "Dip2a", "M2", 2, 0.6666666666666666, 0.327557, 0.982671, 0.0004964331583890222, 0.0004964331583890221, 3, 1, 2.166912666666667, 6.500738, 0.0804092340634467, 0.0804092340634467, 0.7097874123553333, 6.388086711198, 0.00003991781002975899, 0.00003991781002975898),
    @("M2", "Fstl1", "Dip2a", "sCs", 2, 0.6666666666666666, 0.327557, 0.982671, 0.0004964331583890222, 0.0004964331583890221, 3, 1, 4.768749000000001, 14.306247, 0.1769575029162046, 0.1769575029162046, 1.562037116193, 14.058334045737, 0.00008784757207332606, 0.00008784757207332605),
    @("sCs", "Fstl1", "Dip2a", "ECs", 3, 1, 23.31917566666667, 69.957527, 0.03534167191429817, 0.03534167191429816, 3, 1, 9.971009, 29.913027, 0.3700016197528959, 0.3700016197528959, 232.5157104449143, 2092.641394004229, 0.01307647585306575, 0.01307647585306575),
    @("sCs", "Fstl1", "Dip2a", "FAPs", 3, 1, 23.31917566666667, 69.957527, 0.03534167191429817, 0.03534167191429816, 3, 1, 10.04188433333333, 30.125653, 0.3726316432674529, 0.3726316432674529, 234.1684647933479, 2107.516183140131, 0.01316942528124412, 0.01316942528124411),
    @("sCs", "Fstl1", "Dip2a", "M2", 3, 1, 23.31917566666667, 69.957527, 0.03534167191429817, 0.03534167191429816, 3, 1, 2.166912666666667, 6.500738, 0.0804092340634467, 0.0804092340634467, 50.53061712832511, 454.775554154926, 0.002841796769150342, 0.002841796769150342),
    @("sCs", "Fstl1", "Dip2a", "sCs", 3, 1, 23.31917566666667, 69.957527, 0.03534167191429817, 0.03534167191429816, 3, 1, 4.768749000000001, 14.306247, 0.1769575029162046, 0.1769575029162046, 111.203295641241, 1000.829660771169, 0.006253974010837966, 0.006253974010837964)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}